$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths for B (Prices) and C (Costs)
# (target stored widths are 11.42578125 / 10.5703125; the host quantizes
# ColumnWidth to 1/6-character steps, so these inputs land on the closest
# representable stored width.)
$ws.Columns("B").ColumnWidth = 10.6
$ws.Columns("C").ColumnWidth = 9.6

# "Prices" in column B formatted as Euros
$ws.Range("B2:B5").NumberFormat = '"EUR "#,##0.00'

# "Fixed costs" formatted as British Pounds
$ws.Range("B7").NumberFormat = '"GBP "#,##0.00'

# New "Costs" values in column C, formatted as Chinese Yuan
$ws.Range("C2").Value = 170
$ws.Range("C3").Value = 60
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 25
$ws.Range("C2:C5").NumberFormat = '"CNY "#,##0.00'

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update selection
$null = $ws.Range("F6").Select()
